$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Add logging mechanism..." task row entirely - rows below shift up.
$ws.Rows("2:2").Delete()

# The row that shifts into row 2 ("Error handling...") gets a new Estimate value.
$ws.Range("B2").Value = 5

# Update the active selection as recorded in the saved workbook.
$ws.Range("B3").Select()
